$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row values (order chosen so shared-strings table is built in the
#     same order as the target: Tipo Documento, Primer nombre, Segundo Nombre,
#     Primer Apellido, Segundo Apellido, Telefono, Correo Electronico,
#     Numero de documento, Codigo usuario, Jefe, Area, Tipo Usuario) ---
$ws.Range("A1").Value = "Tipo Documento"
$ws.Range("C1").Value = "Primer nombre"
$ws.Range("D1").Value = "Segundo Nombre"
$ws.Range("E1").Value = "Primer Apellido"
$ws.Range("F1").Value = "Segundo Apellido"
$ws.Range("H1").Value = "Telefono"
$ws.Range("G1").Value = "Correo Electronico"
$ws.Range("B1").Value = "Numero de documento"
$ws.Range("K1").Value = "Codigo usuario"
$ws.Range("J1").Value = "Jefe"
$ws.Range("I1").Value = "Area"
$ws.Range("L1").Value = "Tipo Usuario"

# --- Header style: vertical="top" + wrapText="1" applied to A1:N1 (M1 and N1
#     stay valueless but styled, matching the source file). Build the format
#     once on a scratch cell and fan it out with copy/paste-special so the
#     stylesheet ends up with a single extra cellXfs entry, just like the
#     target workbook. ---
$ws.Range("Z1").WrapText = $true
$ws.Range("Z1").VerticalAlignment = -4160
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("A1:N1").PasteSpecial(-4122, -4142, $false, $false) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# --- Column widths (closest attainable values; the runtime snaps
#     ColumnWidth to a 1/6-character grid) ---
$ws.Columns.Item(1).ColumnWidth = 15.5
$ws.Columns.Item(2).ColumnWidth = 25.666666666666668
$ws.Columns.Item(3).ColumnWidth = 18.166666666666668
$ws.Columns.Item(4).ColumnWidth = 18.666666666666668
$ws.Columns.Item(5).ColumnWidth = 16.666666666666668
$ws.Columns.Item(6).ColumnWidth = 17.0
$ws.Columns.Item(7).ColumnWidth = 16.166666666666668
$ws.Columns.Item(8).ColumnWidth = 16.833333333333332
$ws.Columns.Item(9).ColumnWidth = 15.833333333333334
$ws.Columns.Item(11).ColumnWidth = 13.0
$ws.Columns.Item(12).ColumnWidth = 14.333333333333334

# --- Row height, matching the explicit 15pt custom height on the header row ---
$ws.Rows.Item(1).RowHeight = 15

# --- Active selection on J4, as in the target sheet view ---
$ws.Range("J4").Select() | Out-Null
